$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.680.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.695.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.75%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3948'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4023'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.71%  '

$ws.Range("E9").Value = '  +2.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9993'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08765'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.226'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.112'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001316'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.700.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '99.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07046'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.40%  '

$ws.Range("E20").Value = '  +3.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.106'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("E23").Value = '  +3.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.671.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.148'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.337'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.51'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '137.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.207'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.448'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.884.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.84%  '

$ws.Range("E33").Value = '  -3.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08632'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.56%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.139'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2750'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.926'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09146'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.59%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02723'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.44%  '

$ws.Range("E42").Value = '  +1.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7660'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.631'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.86%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7186'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.49%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.223'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9991'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.329'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07992'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.11%  '
